# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right after "总计", built from a copy of
#   the existing "2022-Q2" sheet so it inherits the same header/column
#   formatting (bold+border style on the header row and the "A" index
#   column), then overwrite its data with the two funds reported for
#   2022-Q3.
# - Update the "总计" (summary) roll-up sheet: insert a new top data row for
#   2022-Q3 and append a row for 2021-Q4 so every quarter sheet is listed
#   (the existing 2022-Q2 / 2022-Q1 / 2021-Q4 detail sheets themselves keep
#   their own data unchanged).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (position 2)
#    so fonts/borders/number formats on the header row + index column come
#    along for free, then place it right after "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $wb.Worksheets.Item("总计"))
$q3 = $wb.Worksheets.Item("总计").Next
$q3.Name = "2022-Q3"

# Duplicate the formatting of the existing data row (row 2) down onto row 3
# so the second fund row matches too.
$q3.Range("A2:H2").Copy()
$q3.Range("A3:H3").PasteSpecial(-4122)  # xlPasteFormats

# Row 2: 516530 银华中证现代物流ETF
$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "516530"
$q3.Range("C2").Value = "银华中证现代物流ETF"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "0.89"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "97.53"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "5.06"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0450"
$q3.Range("H2").Value = 3

# Row 3: 516910 富国中证现代物流ETF
$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "516910"
$q3.Range("C3").Value = "富国中证现代物流ETF"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "0.78"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "99.30"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "5.15"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0402"
$q3.Range("H3").Value = 4

# ---------------------------------------------------------------------
# 2. Update "总计": insert the 2022-Q3 row at the top of the data and add
#    the trailing 2021-Q4 row, keeping the per-quarter counters in sync.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the formatting already used by the data rows (incl. the
# bold/bordered style on the "A" index column) down onto the new row 5.
$total.Range("A4:D4").Copy()
$total.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.09

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.03

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.02

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.02
